# comenzando pruebas con lexer/parser dinamico que permita anidamientos condicionales
#
# Content changes applied to Hoja1:
#  - B2/B3/B4 measured pressure values corrected from 9.2 to 2.1
#  - B6 value corrected from 4 to 0
#  - a new data row (7) is appended for PM_IPA_CENTRIFUGADO_MARCHA, re-using the
#    same "fill"/"=" tagging columns as the existing rows and a new dynamic
#    condition (OR instead of AND) to flag the unstable-pressure case
#  - a stray formatted (underlined) cell A11 is left below the table, matching
#    the in-progress edit of the new parser/lexer nesting experiment
#  - the view is refreshed (zoom + active cell) and the print setup is
#    normalised to a single page

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- fix the previously wrong pressure readings in column B -------------
$ws.Range("B2").Value = 2.1
$ws.Range("B3").Value = 2.1
$ws.Range("B4").Value = 2.1
$ws.Range("B6").Value = 0

# --- append row 7: PM_IPA_CENTRIFUGADO_MARCHA with an OR() condition ----
# Clone row 6's layout/formatting/formula first so the new row inherits the
# same styles (centered number columns, "=" / "fill" tag cells), then patch
# in the row-specific values and formula.
$ws.Range("A6:F6").Copy()
$ws.Range("A7:F7").PasteSpecial()

$ws.Cells.Item(7, 1).Value = "PM_IPA_CENTRIFUGADO_MARCHA"
$ws.Cells.Item(7, 2).Value = 0
$ws.Cells.Item(7, 3).Value = 1
$ws.Cells.Item(7, 4).Formula = '=IF(OR(B7<3,B7>7),"presion no estable","presion estable")'

# --- stray underlined, empty cell a couple of rows below the table ------
$ws.Cells.Item(11, 1).Font.Underline = $true

# --- refresh the view: zoom to 100% and move the active selection -------
$excel.ActiveWindow.Zoom = 100
[void]$ws.Range("A9").Select()

# --- normalise the print setup (100% scale, single page) ----------------
$ps = $ws.PageSetup
$ps.Zoom = 100
